# Add a new "2021" column (S) to the table, mirroring the style of column R
# in each row, and update the selection to match the post-edit state
# recorded in the target workbook (Q19, single cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4): new year value
$ws.Range("S4").Value = 2021

# Data rows 5-14: new values for the "2021" column
$values = @{
    5  = 6.1
    6  = 1.6
    7  = 3.6
    8  = 27.2
    9  = 7.2
    10 = 2.6
    11 = 12.5
    12 = 6.4
    13 = 5.2
    14 = 0.9
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# Copy formatting from column R so column S visually matches (borders, etc.)
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the recorded post-edit state
$ws.Range("Q19").Select()
